$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) to bold
$ws.Range("A1:C1").Font.Bold = $true

# Row 4: add "Not O(1) get information or data" in A4 and "LookUp" in C4 (B4 "Find" stays)
$ws.Range("A4").Value = "Not O(1) get information or data"
$ws.Range("C4").Value = "LookUp"

# Remove old "Begin"/"End" content from rows 13 and 14 (without shifting rows 15/16)
$ws.Range("B13").ClearContents()
$ws.Range("B14").ClearContents()

# Page setup: portrait orientation (matches committed pageSetup element)
$ws.PageSetup.Orientation = 1

# Set the selection as shown in diff
$ws.Range("C18").Select()

$wb.Save()
